$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New headers in row 1: X1 = "PriceChange", Y1 = "UpDown"
# ---------------------------------------------------------------------
$ws.Cells.Item(1,24).Value = "PriceChange"
$ws.Cells.Item(1,25).Value = "UpDown"

# ---------------------------------------------------------------------
# Row 2 (existing row) updates
# ---------------------------------------------------------------------
$ws.Cells.Item(2,1).Value  = 42633.878912037035   # A2 Date
$ws.Cells.Item(2,2).Value  = -12                  # B2 ScoreFinal
$ws.Cells.Item(2,3).Value  = "Buy"                # C2 Verdict
$ws.Cells.Item(2,16).Value = "Random"             # P2 Method (same text, kept)
$ws.Cells.Item(2,18).Value = 1.76                 # R2 PEG
$ws.Cells.Item(2,19).Value = 0.111                # S2 200Moving%
$ws.Cells.Item(2,20).Value = 5.45                 # T2 50Moving%
$ws.Cells.Item(2,21).Value = 4.84                 # U2 PriceBook

# New trailing cells on row 2
$ws.Cells.Item(2,24).Value = -1.6100000000000136  # X2 PriceChange
$ws.Cells.Item(2,25).Value = "Down"               # Y2 UpDown

# ---------------------------------------------------------------------
# Row 3 (brand new row) - copy formatting from row 2 first so number
# formats (date / percentage) line up with the existing style indexes,
# then overwrite with the correct values.
# ---------------------------------------------------------------------
$ws.Range("A2:W2").Copy($ws.Range("A3:W3"))

$ws.Cells.Item(3,1).Value  = 42633.880312499998   # A3 Date
$ws.Cells.Item(3,2).Value  = 0                    # B3 ScoreFinal
$ws.Cells.Item(3,3).Value  = "Buy"                # C3 Verdict
$ws.Cells.Item(3,4).Value  = 0                    # D3 totalSentiment
$ws.Cells.Item(3,5).Value  = 0                    # E3 wordCount
$ws.Cells.Item(3,6).Value  = 0                    # F3 sentenceCount
$ws.Cells.Item(3,7).Value  = 0                    # G3 posWordPercentage
$ws.Cells.Item(3,8).Value  = 0                    # H3 negWordPercentage
$ws.Cells.Item(3,9).Value  = 0                    # I3 posPhrasePercentage
$ws.Cells.Item(3,10).Value = 0                    # J3 negPhrasePercentage
$ws.Cells.Item(3,11).Value = 0                    # K3 ElapsedMs
$ws.Cells.Item(3,12).Value = 0                    # L3 posWordCount
$ws.Cells.Item(3,13).Value = 0                    # M3 negWordCount
$ws.Cells.Item(3,14).Value = 0                    # N3 positivePhraseCount
$ws.Cells.Item(3,15).Value = 0                    # O3 negativePhraseCount
$ws.Cells.Item(3,16).Value = "Random"             # P3 Method
$ws.Cells.Item(3,17).Value = 0                    # Q3 RSI
$ws.Cells.Item(3,18).Value = 1.76                 # R3 PEG
$ws.Cells.Item(3,19).Value = 0.111                # S3 200Moving%
$ws.Cells.Item(3,20).Value = 5.45                 # T3 50Moving%
$ws.Cells.Item(3,21).Value = 4.84                 # U3 PriceBook
$ws.Cells.Item(3,22).Value = 2.2799999999999998   # V3 Dividend
$ws.Cells.Item(3,23).Value = 0                    # W3 Bollinger
